$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells - no ambiguous auto-conversion, set directly.
$ws.Range("A2").Value = "Daniel "
$ws.Range("C2").Value = "(917) 975-2625"
$ws.Range("D2").Value = "dshifrin5@gmail.com"
$ws.Range("F2").Value = "www"
$ws.Range("G2").Value = "2025-05-27 21-43-00"

# B2/E2 look like a date / number respectively, so Excel would otherwise
# auto-convert them on assignment. Force text entry via NumberFormat "@",
# then restore the default "Normal" style so no residual formatting sticks.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2002-04-25"
$ws.Range("B2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1181"
$ws.Range("E2").Style = "Normal"
